# fix Dob loader issue + Address Teacher Issue
#
# - "Học sinh" (Students) sheet: the "Ngày sinh" (Dob) column is column B.
# - "Giáo viên" (Teachers) sheet: the "Ngày sinh" (Dob) column is column D.
#
# The Dob columns (and, as a side effect of the shared cell style used for
# the phone-number columns, the phone-number columns too) get their font
# normalized to an explicit black, and get left-aligned so the Dob text
# values display consistently instead of using the sheet's general
# formatting. The header/data rows on the Students sheet also grow
# slightly taller to match.

$wb = $excel.ActiveWorkbook

$students = $wb.Worksheets.Item("Học sinh")
$teachers = $wb.Worksheets.Item("Giáo viên")

# --- Students ("Học sinh") sheet -------------------------------------------------

# Row heights: header row grows from 21 -> 20.25, data rows from 18.75 -> 19.5
$students.Rows.Item(1).RowHeight = 20.25
$students.Range("2:61").RowHeight = 19.5

# Dob column (B) - left align + normalize font color to black
$students.Range("B2:B61").HorizontalAlignment = -4131
$students.Range("B2:B61").Font.Color = 0

# Phone number columns that shared the old Dob/box style (E, H, L) pick up
# the same font color normalization
$students.Range("E2:E61").Font.Color = 0
$students.Range("H2:H61").Font.Color = 0
$students.Range("L2:L61").Font.Color = 0

# --- Teachers ("Giáo viên") sheet -------------------------------------------------

# Dob column (D) - left align + normalize font color to black
$teachers.Range("D2:D21").HorizontalAlignment = -4131
$teachers.Range("D2:D21").Font.Color = 0

# Phone number columns (B, G) pick up the same font color normalization
$teachers.Range("B2:B21").Font.Color = 0
$teachers.Range("G2:G21").Font.Color = 0
